$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = 50
$ws.Range("L6").Value = 150
$ws.Range("N6").Value = -374

$ws.Range("H15").Value = 2209.08
$ws.Range("I15").Value = 2209.08
$ws.Range("K15").Value = 6627.24
$ws.Range("M15").Value = -6458.24

$ws.Range("H29").Value = 99
$ws.Range("I29").Value = 99
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 297
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -16
$ws.Range("N29").ClearContents()

$ws.Range("H38").Value = 1220.5
$ws.Range("I38").Value = 162.8
$ws.Range("J38").Value = 6509
$ws.Range("K38").Value = 488.4
$ws.Range("L38").Value = 19527
$ws.Range("M38").Value = -116.4
$ws.Range("N38").Value = -20271

$ws.Range("H58").Value = 6383.5
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 12017
$ws.Range("K58").Value = 2250
$ws.Range("L58").Value = 36051
$ws.Range("M58").Value = -2100
$ws.Range("N58").Value = -36351

$ws.Range("H70").Value = 750
$ws.Range("I70").Value = 750
$ws.Range("K70").Value = 2250
$ws.Range("M70").Value = -1980

$ws.Range("H73").Value = 750
$ws.Range("I73").Value = 750
$ws.Range("K73").Value = 2250
$ws.Range("M73").Value = -1314

$ws.Range("H97").Value = 399.5
$ws.Range("J97").Value = 399.5
$ws.Range("L97").Value = 1198.5
$ws.Range("N97").Value = -2190.5

$ws.Range("H137").Value = 2159.95
$ws.Range("I137").Value = 992.75
$ws.Range("J137").Value = 2938.0833
$ws.Range("K137").Value = 2978.25
$ws.Range("L137").Value = 8814.249899999999
$ws.Range("M137").Value = -428.25
$ws.Range("N137").Value = -13914.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5516.25
$ws.Range("I61").Value = 5828.143
$ws.Range("K61").Value = 5828.143
$ws.Range("M61").Value = -5616.143

$ws.Range("H122").Value = 1681.2
$ws.Range("I122").Value = 1681.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5043.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2593.6
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2924.9473
$ws.Range("I132").Value = 2359.8462
$ws.Range("K132").Value = 7079.5386
$ws.Range("M132").Value = -4549.5386

$ws.Range("H136").Value = 5516.25
$ws.Range("I136").Value = 5828.143
$ws.Range("K136").Value = 17484.429
$ws.Range("M136").Value = -14934.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 627.3
$ws.Range("J80").Value = 312
$ws.Range("L80").Value = 312
$ws.Range("N80").Value = -2308

$ws.Range("H83").Value = 627.3
$ws.Range("J83").Value = 312
$ws.Range("L83").Value = 1560
$ws.Range("N83").Value = -11544

$ws.Range("H105").Value = 3098.5
$ws.Range("I105").Value = 3098.5
$ws.Range("K105").Value = 3098.5
$ws.Range("M105").Value = -1351.5

$ws.Range("H134").Value = 3351.1
$ws.Range("I134").Value = 2941.2354
$ws.Range("K134").Value = 8823.706200000001
$ws.Range("M134").Value = -6288.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9999.5
$ws.Range("I16").Value = 9999.5
$ws.Range("K16").Value = 9999.5
$ws.Range("M16").Value = -9712.5

$ws.Range("H43").Value = 19999.5
$ws.Range("J43").Value = 19999.5
$ws.Range("L43").Value = 19999.5
$ws.Range("N43").Value = -20367.5

$ws.Range("H58").Value = 2144.5454
$ws.Range("I58").Value = 2065.889
$ws.Range("K58").Value = 2065.889
$ws.Range("M58").Value = -1862.889

$ws.Range("H62").Value = 3283.4285
$ws.Range("I62").Value = 3077.8
$ws.Range("J62").Value = 3797.5
$ws.Range("K62").Value = 3077.8
$ws.Range("L62").Value = 3797.5
$ws.Range("M62").Value = -2453.8
$ws.Range("N62").Value = -5045.5

$ws.Range("H65").Value = 3283.4285
$ws.Range("I65").Value = 3077.8
$ws.Range("J65").Value = 3797.5
$ws.Range("K65").Value = 15389
$ws.Range("L65").Value = 18987.5
$ws.Range("M65").Value = -12269
$ws.Range("N65").Value = -25227.5

$ws.Range("H101").Value = 19999.5
$ws.Range("J101").Value = 19999.5
$ws.Range("L101").Value = 19999.5
$ws.Range("N101").Value = -26489.5

$ws.Range("H113").Value = 9999.5
$ws.Range("I113").Value = 9999.5
$ws.Range("K113").Value = 9999.5
$ws.Range("M113").Value = -7829.5

$ws.Range("H136").Value = 2144.5454
$ws.Range("I136").Value = 2065.889
$ws.Range("K136").Value = 6197.667
$ws.Range("M136").Value = -3647.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 750
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 1500
$ws.Range("M17").Value = -1331

$ws.Range("H37").Value = 186363
$ws.Range("J37").Value = 186363
$ws.Range("L37").Value = 559089
$ws.Range("N37").Value = -559313

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.5
$ws.Range("I2").Value = 106.666664
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 106.666664
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 6.333336000000003
$ws.Range("N2").Value = -276

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H43").Value = 15831.667
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7849

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H46").Value = 34920
$ws.Range("J46").Value = 34920
$ws.Range("L46").Value = 34920
$ws.Range("N46").Value = -35232

$ws.Range("H80").Value = 6374.5
$ws.Range("I80").Value = 6374.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 6374.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -5376.5
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 6374.5
$ws.Range("I83").Value = 6374.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 31872.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -26880.5
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2414.1667
$ws.Range("I46").Value = 897.5
$ws.Range("J46").Value = 5447.5
$ws.Range("K46").Value = 897.5
$ws.Range("L46").Value = 5447.5
$ws.Range("M46").Value = -709.5
$ws.Range("N46").Value = -5823.5

$ws.Range("H82").Value = 1524.75
$ws.Range("I82").Value = 2100
$ws.Range("J82").Value = 949.5
$ws.Range("K82").Value = 2100
$ws.Range("L82").Value = 949.5
$ws.Range("M82").Value = -1739
$ws.Range("N82").Value = -1671.5

$ws.Range("H85").Value = 1524.75
$ws.Range("I85").Value = 2100
$ws.Range("J85").Value = 949.5
$ws.Range("K85").Value = 2100
$ws.Range("L85").Value = 949.5
$ws.Range("M85").Value = -852
$ws.Range("N85").Value = -3445.5

$ws.Range("H131").Value = 47500
$ws.Range("J131").Value = 47500
$ws.Range("L131").Value = 47500
$ws.Range("N131").Value = -57580

$ws.Range("H136").Value = 3858.4285
$ws.Range("I136").Value = 3858.4285
$ws.Range("K136").Value = 11575.2855
$ws.Range("M136").Value = -9025.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 320
$ws.Range("I96").Value = 320
$ws.Range("K96").Value = 320
$ws.Range("M96").Value = 1053

$ws.Range("H107").Value = 1001
$ws.Range("J107").Value = 1001
$ws.Range("L107").Value = 3003
$ws.Range("N107").Value = -6843

$ws.Range("H132").Value = 1285.64
$ws.Range("I132").Value = 855.2105
$ws.Range("K132").Value = 2565.6315
$ws.Range("M132").Value = -35.63149999999996

$ws.Range("H136").Value = 2628
$ws.Range("I136").Value = 1759.8
$ws.Range("J136").Value = 6969
$ws.Range("K136").Value = 5279.4
$ws.Range("L136").Value = 20907
$ws.Range("M136").Value = -2729.4
$ws.Range("N136").Value = -26007
